$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to append (231-235), mirroring the existing layout/format of row 230.
# Columns: A Color, B Comment, C Document group, D Document name, E Code,
#          F Begin, G End, H Weight score, I Segment, J Area, K Coverage %,
#          L Author, M Creation date

$rowsData = @(
    @{ Row=231; A="●"; B=""; C=""; D="2697"; E="Bacteria:Binomial (genus species)"; F="1: 1003"; G="1: 1023"; H=0; I="Staphylococcus aureus"; J=21; K=0.080912; L="Sonia"; M="11/8/18 14:45:00" },
    @{ Row=232; A="●"; B=""; C=""; D="3651"; E="Bacteria:Binomial (genus species)"; F="1: 3451"; G="1: 3471"; H=0; I="Neisseria gonorrhoeae"; J=21; K=0.046607; L="Sonia"; M="11/8/18 14:46:00" },
    @{ Row=233; A="●"; B=""; C=""; D="3910"; E="Bacteria:Binomial (genus species)"; F="1: 1015"; G="1: 1035"; H=0; I="Mycoplasma genitalium"; J=21; K=0.107714; L="Sonia"; M="11/8/18 14:47:00" },
    @{ Row=234; A="●"; B=""; C=""; D="5251"; E="Bacteria:Binomial (genus species)"; F="1: 2652"; G="1: 2662"; H=0; I="S. enterica"; J=11; K=0.062825; L="Sonia"; M="11/8/18 14:48:00" },
    @{ Row=235; A="●"; B=""; C=""; D="5251"; E="Bacteria:Strain"; F="1: 2673"; G="1: 2683"; H=0; I="typhimurium"; J=11; K=0.062825; L="Sonia"; M="11/8/18 14:48:00" }
)

foreach ($rd in $rowsData) {
    $r = $rd.Row

    # Set text-like values first (as text, via leading apostrophe so numeric-looking
    # strings like document names stay text) - style gets fixed up right after.
    $ws.Cells.Item($r, 1).Value = "'" + $rd.A
    $ws.Cells.Item($r, 2).Value = "'" + $rd.B
    $ws.Cells.Item($r, 3).Value = "'" + $rd.C
    $ws.Cells.Item($r, 4).Value = "'" + $rd.D
    $ws.Cells.Item($r, 5).Value = "'" + $rd.E
    $ws.Cells.Item($r, 6).Value = "'" + $rd.F
    $ws.Cells.Item($r, 7).Value = "'" + $rd.G
    $ws.Cells.Item($r, 8).Value = $rd.H
    $ws.Cells.Item($r, 9).Value = "'" + $rd.I
    $ws.Cells.Item($r, 10).Value = $rd.J
    $ws.Cells.Item($r, 11).Value = $rd.K
    $ws.Cells.Item($r, 12).Value = "'" + $rd.L
    $ws.Cells.Item($r, 13).Value = "'" + $rd.M

    # Copy the formatting (styles) from row 230, the last existing data row,
    # onto the newly populated row, without touching values.
    $ws.Range("A230:M230").Copy()
    $ws.Range("A" + $r + ":M" + $r).PasteSpecial(-4122)

    $ws.Rows.Item($r).RowHeight = 16
}

$excel.CutCopyMode = 0
